$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 27.83096748370949
$ws.Range("C2").Value = 28.84720169803776
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 27.79020060783168
$ws.Range("C3").Value = 28.81483507752534
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 27.8574254088211
$ws.Range("C4").Value = 28.84831240086282
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 27.84740778959875
$ws.Range("C5").Value = 28.88794157735925
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 27.84937074870759
$ws.Range("C6").Value = 28.85522678685512
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 27.81816264608317
$ws.Range("C7").Value = 28.84828755242988
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 27.8390060300491
$ws.Range("C8").Value = 28.84712767239659
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 27.82131605091118
$ws.Range("C9").Value = 28.84829255545864
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 27.81742838227542
$ws.Range("C10").Value = 28.83940459846023
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 27.83695446924043
$ws.Range("C11").Value = 28.83634974466139
$ws.Range("B12").Value = 27.83082396172279
$ws.Range("C12").Value = 28.8472979664047

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 19.81360772831505
$ws.Range("C2").Value = 26.1509069501992
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 19.85848040293938
$ws.Range("C3").Value = 26.1855511587085
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 19.82127581465778
$ws.Range("C4").Value = 26.14690495315903
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 19.85354781174944
$ws.Range("C5").Value = 26.19741383258339
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 19.82843320122187
$ws.Range("C6").Value = 26.14554304271921
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 19.83117714871628
$ws.Range("C7").Value = 26.14261871513665
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 19.8267338368039
$ws.Range("C8").Value = 26.13777282359939
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 19.81918576695703
$ws.Range("C9").Value = 26.14208799493854
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 19.856709838037
$ws.Range("C10").Value = 26.21422641470159
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 19.83789688697364
$ws.Range("C11").Value = 26.13966739550372
$ws.Range("B12").Value = 19.83470484363714
$ws.Range("C12").Value = 26.16026932812492

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 14.961609172719
$ws.Range("C2").Value = 21.11189937134795
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 14.9713443928894
$ws.Range("C3").Value = 21.15951761530154
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 14.97235008889229
$ws.Range("C4").Value = 21.1558129446036
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 14.9919493135604
$ws.Range("C5").Value = 21.19527039322807
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 14.96577627950721
$ws.Range("C6").Value = 21.17988395771821
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 14.98428047187979
$ws.Range("C7").Value = 21.17596760470006
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 14.97449471337259
$ws.Range("C8").Value = 21.16236955987747
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 14.95685614272055
$ws.Range("C9").Value = 21.11077210605161
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 14.97950404299135
$ws.Range("C10").Value = 21.15126370954349
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 14.96678275407911
$ws.Range("C11").Value = 21.1296588540973
$ws.Range("B12").Value = 14.97249473726117
$ws.Range("C12").Value = 21.15324161164693
